# Fix the misspelled "Brookyn Nets" entry to "Brooklyn Nets" and move the
# active selection to the corrected cell (B3), matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Brooklyn Nets"
$ws.Range("B3").Select() | Out-Null
